$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 11:02:23"
$wsZhCn.Range("H2").Value = "2016-03-21 11:02:41"

# de-de sheet: update Correspond Handoff/Handback Datetime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 11:02:26"
$wsDeDe.Range("H2").Value = "2016-03-21 11:02:46"
